$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$row = 69

# New data row appended below the existing table (row 68 was the previous last row)
$ws.Cells.Item($row, 1).Value = "Emna"
$ws.Cells.Item($row, 2).Value = "Boughariou"
$ws.Cells.Item($row, 3).Value = "Université de Sfax"
$ws.Cells.Item($row, 4).Value = "Tunisie"
$ws.Cells.Item($row, 5).Value = "wdvUCRwAAAAJ"
$ws.Cells.Item($row, 6).Value = "F"
$ws.Cells.Item($row, 7).Value = 1992
$ws.Cells.Item($row, 8).Value = "Informatique, Mathématiques et Ingénierie"

# Column F (Genre) carries a small Arial font elsewhere in the sheet - copy that
# formatting onto the new cell instead of rebuilding it by hand.
$ws.Cells.Item(68, 6).Copy() | Out-Null
$ws.Cells.Item($row, 6).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reposition the view / selection as it ended up after the edit
$ws.Range("E71").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 51
$excel.ActiveWindow.ScrollColumn = 4
